$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before the existing column B (gold/gem data moves right).
$ws.Range("B:E").EntireColumn.Insert()

# New header row (row 2)
$ws.Range("B2").Value = "noneGradeRate"
$ws.Range("C2").Value = "bronzeGradeRate"
$ws.Range("D2").Value = "silverGradeRate"
$ws.Range("E2").Value = "goldGradeRate"

# New type row (row 3)
$ws.Range("B3").Value = "float"
$ws.Range("C3").Value = "float"
$ws.Range("D3").Value = "float"
$ws.Range("E3").Value = "float"

# New data rows (4 through 17): constant grade-rate values per column.
for ($r = 4; $r -le 17; $r++) {
    $ws.Cells.Item($r, 2).Value = 4
    $ws.Cells.Item($r, 3).Value = 3
    $ws.Cells.Item($r, 4).Value = 2
    $ws.Cells.Item($r, 5).Value = 1
}

# Column widths to fit the new header text (bestFit columns).
$ws.Columns.Item(2).ColumnWidth = 14.28515625
$ws.Columns.Item(3).ColumnWidth = 15.7109375
$ws.Columns.Item(4).ColumnWidth = 14.28515625
$ws.Columns.Item(5).ColumnWidth = 13.7109375
$ws.Columns.Item(6).ColumnWidth = 11.5703125

# Selection moved in the source file.
$ws.Range("I13").Select()
